$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update invoice date/time (I4)
$ws.Range("I4").Value = "2023-11-29 17:46:50"

# 2. Add a customer row: "Khách:" label (H5) and hashed customer value (I5).
#    Reuse the existing formats from the "Ngày:" row (H4/I4/L4) via format-only
#    paste so the new cells share style indices instead of creating new ones.
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("I4").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("H5").Value = "Khách:"
$ws.Range("I5").Value = "w6Cfy19ejrOSkW3NkazDfOeduDNzhrmgTdpglJOLc9M="

$ws.Range("I5:L5").Merge()
$ws.Range("J5").ClearFormats()
$ws.Range("K5").ClearFormats()
$ws.Range("L4").Copy()
$ws.Range("L5").PasteSpecial(-4122)

# 3. Update the order line (row 11): new item, quantity, price; drop the discount
$ws.Range("B11").Value = "Gà nướng muối ớt (1 con)"
$ws.Range("F11").Value = 1
$ws.Range("K11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G11").Value = "130.000 ₫"
$ws.Range("I11").Value = ""
$ws.Range("K11").Value = "130.000 ₫"

# 4. Remove the second order line entirely (old row 12); rows below shift up by one
$ws.Rows.Item(12).Delete()

# 5. Update totals (post-shift: old rows 16/17/18 are now 15/16/17)
$ws.Range("I15").Value = "130.000 ₫"
$ws.Range("I16").Value = "200.000 ₫"
$ws.Range("I17").Value = "-70.000 ₫"

Write-Host "done"
